$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: make a cell a "no data" text placeholder that reuses the existing
# N/A-style shared strings ("0" or "***.*"), matching the style of a
# neighboring cell that already uses that placeholder style (style index 14).
# We first write the value as plain text (so the shared string resolves to
# the pre-existing "0"/"***.*" entry), then paste just the *format* from a
# reference cell that already carries the desired style, which overwrites
# the quote-prefix styling introduced by the text assignment without
# touching the value/type we just set.
# ---------------------------------------------------------------------------
function Set-NAPlaceholder {
    param($targetAddr, $text, $formatSourceAddr)
    $ws.Range($targetAddr).Value = "'" + $text
    $ws.Range($formatSourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# Helper: turn a "no data" text placeholder cell into a real number, copying
# the numeric style (format) from a neighboring numeric cell that already
# has the desired style.
# ---------------------------------------------------------------------------
function Set-NumberWithStyle {
    param($targetAddr, $value, $formatSourceAddr)
    $ws.Range($targetAddr).Value = $value
    $ws.Range($formatSourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null
}

# ===========================================================================
# Title block text updates (report volume/number + week-of dates)
# ===========================================================================
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# ===========================================================================
# Row 15 - Rape
# ===========================================================================
Set-NumberWithStyle "D15" 1 "G15"
Set-NumberWithStyle "E15" -100 "H15"
$ws.Range("G15").Value = 3
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -40
$ws.Range("L15").Value = -40

# ===========================================================================
# Row 16 - Robbery
# ===========================================================================
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -4.761904761904
$ws.Range("I16").Value = 104
$ws.Range("J16").Value = 129
$ws.Range("K16").Value = -19.379844961240
$ws.Range("L16").Value = -29.729729729729
$ws.Range("M16").Value = -30.201342281879

# ===========================================================================
# Row 17 - Fel. Assault
# ===========================================================================
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -9.523809523809
$ws.Range("I17").Value = 158
$ws.Range("J17").Value = 187
$ws.Range("K17").Value = -15.508021390374
$ws.Range("L17").Value = -16.402116402116
$ws.Range("M17").Value = 31.666666666666

# ===========================================================================
# Row 18 - Burglary
# ===========================================================================
Set-NAPlaceholder "C18" "0" "C15"
Set-NAPlaceholder "D18" "0" "D15"
Set-NAPlaceholder "E18" "***.*" "E15"
$ws.Range("L18").Value = -19.354838709677
$ws.Range("M18").Value = 25

# ===========================================================================
# Row 19 - Gr. Larceny
# ===========================================================================
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 700
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 74.074074074074
$ws.Range("I19").Value = 229
$ws.Range("J19").Value = 198
$ws.Range("K19").Value = 15.656565656565
$ws.Range("L19").Value = 16.243654822335
$ws.Range("M19").Value = 37.125748502994

# ===========================================================================
# Row 20 - G.L.A.
# ===========================================================================
$ws.Range("C20").Value = 2
Set-NAPlaceholder "D20" "0" "D15"
Set-NAPlaceholder "E20" "***.*" "E15"
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 77.777777777777
$ws.Range("I20").Value = 99
$ws.Range("K20").Value = 5.319148936170
$ws.Range("L20").Value = 147.5
$ws.Range("M20").Value = 135.714285714286

# ===========================================================================
# Row 21 - TOTAL
# ===========================================================================
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 76.470588235294
$ws.Range("F21").Value = 105
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = 20.689655172413
$ws.Range("I21").Value = 675
$ws.Range("J21").Value = 682
$ws.Range("K21").Value = -1.026392961876
$ws.Range("L21").Value = -0.589101620029
$ws.Range("M21").Value = 22.504537205081

# ===========================================================================
# Row 22 - Transit
# ===========================================================================
$ws.Range("C22").Value = 2
Set-NAPlaceholder "D22" "0" "D15"
Set-NAPlaceholder "E22" "***.*" "E15"
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 11
$ws.Range("K22").Value = -45
$ws.Range("L22").Value = 22.222222222222
$ws.Range("M22").Value = 37.5

# ===========================================================================
# Row 23 - Housing
# ===========================================================================
$ws.Range("L23").Value = -42.857142857142

# ===========================================================================
# Row 24 - Petit Larceny
# ===========================================================================
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -33.333333333333
$ws.Range("F24").Value = 86
$ws.Range("G24").Value = 86
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 617
$ws.Range("J24").Value = 933
$ws.Range("K24").Value = -33.869239013933
$ws.Range("L24").Value = -34.010695187165
$ws.Range("M24").Value = 73.314606741573

# ===========================================================================
# Row 25 - Misd. Assault
# ===========================================================================
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 36
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 5.882352941176
$ws.Range("I25").Value = 258
$ws.Range("J25").Value = 268
$ws.Range("K25").Value = -3.731343283582
$ws.Range("L25").Value = -0.386100386100
$ws.Range("M25").Value = -18.095238095238

# ===========================================================================
# Row 26 - UCR Rape*
# ===========================================================================
Set-NumberWithStyle "D26" 1 "G26"
Set-NumberWithStyle "E26" -100 "H26"
$ws.Range("G26").Value = 3
$ws.Range("J26").Value = 16
$ws.Range("K26").Value = -50
$ws.Range("L26").Value = -50

# ===========================================================================
# Row 27 - Other Sex Crimes
# ===========================================================================
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 26
$ws.Range("J27").Value = 36
$ws.Range("K27").Value = -27.777777777777
$ws.Range("L27").Value = -13.333333333333

# ===========================================================================
# Row 30 - Hate Crimes
# ===========================================================================
Set-NAPlaceholder "D30" "0" "D15"
Set-NAPlaceholder "E30" "***.*" "E15"
